$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new values for columns C (nombre_aides) and E (montant_total)
$updates = @{
    4   = @{ C = 46789;  E = 126189040 }
    8   = @{ C = 181369; E = 653049910 }
    115 = @{ C = 17560;  E = 38644994 }
    152 = @{ C = 126053; E = 716083919 }
    164 = @{ C = 50585;  E = 168943069 }
    168 = @{ C = 285114; E = 1213524897 }
    169 = @{ C = 562672; E = 1286084501 }
    170 = @{ C = 367559; E = 2848069929 }
    171 = @{ C = 115224; E = 448839617 }
    174 = @{ C = 357363; E = 1019961276 }
    175 = @{ C = 125687; E = 815654401 }
    179 = @{ C = 235796; E = 813626581 }
    180 = @{ C = 141529; E = 341235629 }
    220 = @{ C = 4714;   E = 11992288 }
    237 = @{ C = 58310;  E = 172543972 }
    239 = @{ C = 84901;  E = 500391684 }
    255 = @{ C = 141372; E = 414550820 }
    303 = @{ C = 40039;  E = 131830133 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("E$row").Value = $vals.E
}
